# Replace the four "Kampagnendaten 2018 für das Sternbild Perseus: ..." campaign-date
# paragraphs with the new Taurus campaign date text, collapsing every run in
# each paragraph into a single, unformatted run.

$d = $word.ActiveDocument

$oldText = "Kampagnendaten 2018 für das Sternbild Perseus: 30. Oktober - 8. November und 29. November - 8. Dezember"
$newText = "Kampagnendaten Taurus: 16. bis 25. Januar"

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $r = $para.Range
    # Exclude the trailing paragraph mark from the range we inspect/replace.
    $r.End = $r.End - 1
    if ($r.Text -eq $oldText) {
        # First clear all existing runs (and their formatting) ...
        $r.Text = ""
        # ... then type the replacement into what is now an empty paragraph,
        # so the new run picks up no leftover run-level formatting.
        $r2 = $para.Range
        $r2.End = $r2.End - 1
        $r2.Text = $newText
    }
}
